$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - copy the format of the existing header row (G1)
# so it picks up the same bold/border/centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values under the Save header.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
